$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the end of the data block (after row 16) to extend to 19 rows
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(17).Insert()

# Copy formatting (border/font/alignment) from the last original data row (16) to the new rows (17:19)
$ws.Range("A16:M16").Copy()
$ws.Range("A17:M17").PasteSpecial(-4122)
$ws.Range("A18:M18").PasteSpecial(-4122)
$ws.Range("A19:M19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rewrite rows 10..19 (A=8..17) with final values (labels + intensities)
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.9528535474453399
$ws.Range("D10").Value = 1.1851820878893
$ws.Range("E10").Value = 0.9489456796225956
$ws.Range("F10").Value = 0.9528535474453399
$ws.Range("G10").Value = 1.102144700832245
$ws.Range("H10").Value = 0.8620159336131725
$ws.Range("I10").Value = 0.9450172605927911
$ws.Range("J10").Value = 1.1851820878893
$ws.Range("K10").Value = 1.067063883755948
$ws.Range("L10").Value = 1.009958715600644
$ws.Range("M10").Value = 0.9993598683325741

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9992396605434293
$ws.Range("D11").Value = 0.8910338054020623
$ws.Range("E11").Value = 1.030571606813469
$ws.Range("F11").Value = 0.9992396605434293
$ws.Range("G11").Value = 0.9317422398444506
$ws.Range("H11").Value = 1.095468658029462
$ws.Range("I11").Value = 1.023772164988507
$ws.Range("J11").Value = 0.8910338054020623
$ws.Range("K11").Value = 0.9608027061077656
$ws.Range("L11").Value = 0.9800211833255974
$ws.Range("M11").Value = 0.9953046892702301

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9989173332838192
$ws.Range("D12").Value = 0.8915076231967975
$ws.Range("E12").Value = 1.030513198479108
$ws.Range("F12").Value = 0.9989173332838192
$ws.Range("G12").Value = 0.9319470414445374
$ws.Range("H12").Value = 1.095496555503946
$ws.Range("I12").Value = 1.023702305444473
$ws.Range("J12").Value = 0.8915076231967975
$ws.Range("K12").Value = 0.9610104108379527
$ws.Range("L12").Value = 0.979963872060886
$ws.Range("M12").Value = 0.9953473428921135

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9991492332374177
$ws.Range("D13").Value = 0.8912452291395052
$ws.Range("E13").Value = 1.03047600758767
$ws.Range("F13").Value = 0.9991492332374177
$ws.Range("G13").Value = 0.93181969850216
$ws.Range("H13").Value = 1.095388013879673
$ws.Range("I13").Value = 1.023733494529424
$ws.Range("J13").Value = 0.8912452291395052
$ws.Range("K13").Value = 0.9608606183635878
$ws.Range("L13").Value = 0.9800049258005028
$ws.Range("M13").Value = 0.9953019461459752

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.8489320000000012
$ws.Range("D14").Value = 1.452999999999997
$ws.Range("E14").Value = 0.9060360000000003
$ws.Range("F14").Value = 0.8489320000000012
$ws.Range("G14").Value = 1.239003999999999
$ws.Range("H14").Value = 0.7493399999999997
$ws.Range("I14").Value = 0.8847720000000018
$ws.Range("J14").Value = 1.452999999999997
$ws.Range("K14").Value = 1.179517999999999
$ws.Range("L14").Value = 1.014225
$ws.Range("M14").Value = 1.013514

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.76
$ws.Range("D15").Value = 1.83
$ws.Range("E15").Value = 0.8165500000000008
$ws.Range("F15").Value = 0.76
$ws.Range("G15").Value = 1.448349999999997
$ws.Range("H15").Value = 0.5
$ws.Range("I15").Value = 0.79
$ws.Range("J15").Value = 1.83
$ws.Range("K15").Value = 1.323275
$ws.Range("L15").Value = 1.0416375
$ws.Range("M15").Value = 1.02415

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.8566941956096013
$ws.Range("D16").Value = 1.481761670553597
$ws.Range("E16").Value = 0.8903419113472026
$ws.Range("F16").Value = 0.8566941956096013
$ws.Range("G16").Value = 1.257118367232
$ws.Range("H16").Value = 0.7089918062592004
$ws.Range("I16").Value = 0.8767744202752022
$ws.Range("J16").Value = 1.481761670553597
$ws.Range("K16").Value = 1.1860517909504
$ws.Range("L16").Value = 1.02137299328
$ws.Range("M16").Value = 1.011947061879467

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9967212920330168
$ws.Range("D17").Value = 0.9948912011237511
$ws.Range("E17").Value = 0.9950882865136689
$ws.Range("F17").Value = 0.9967212920330168
$ws.Range("G17").Value = 0.9970342577275021
$ws.Range("H17").Value = 0.9936737435347085
$ws.Range("I17").Value = 0.9950740276813855
$ws.Range("J17").Value = 0.9948912011237511
$ws.Range("K17").Value = 0.99498974381871
$ws.Range("L17").Value = 0.9958555179258634
$ws.Range("M17").Value = 0.9954138014356722

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.001704617832466
$ws.Range("D18").Value = 0.9703317057073563
$ws.Range("E18").Value = 1.00197211203712
$ws.Range("F18").Value = 1.001704617832466
$ws.Range("G18").Value = 0.9820684144004888
$ws.Range("H18").Value = 1.01217316590617
$ws.Range("I18").Value = 1.002125710800175
$ws.Range("J18").Value = 0.9703317057073563
$ws.Range("K18").Value = 0.9861519088722384
$ws.Range("L18").Value = 0.9939282633523518
$ws.Range("M18").Value = 0.9950626211139627

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.014913140328588
$ws.Range("D19").Value = 0.909582521458742
$ws.Range("E19").Value = 1.015629259677803
$ws.Range("F19").Value = 1.014913140328588
$ws.Range("G19").Value = 0.945077579816057
$ws.Range("H19").Value = 1.054597876214529
$ws.Range("I19").Value = 1.018177930413418
$ws.Range("J19").Value = 0.909582521458742
$ws.Range("K19").Value = 0.9626058905682726
$ws.Range("L19").Value = 0.9887595154484301
$ws.Range("M19").Value = 0.9929963846515228
